$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point: insert a row at row 7 (pushing the existing
# rows 7-21 down to 8-22) and populate it with the latest observation.
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44575
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112030
$ws.Range("G7").Value = "Poroto granado"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 35000
$ws.Range("L7").Value = 35000
$ws.Range("M7").Value = 35000
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 1400
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
